$d = $word.ActiveDocument

# --- 1. Remove the stray "00000000000000000" run from the first
#        paragraph ("Parte 1...") and drop the _GoBack bookmark that
#        currently sits there (it will be relocated below). ---
$d.Content.Find.Execute("00000000000000000", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2. Give the "CONFLICTO DE VERSIONES" paragraph an explicit paragraph
#        mark formatting (bold, size 24) and append a new, empty paragraph
#        after it that carries the (relocated) _GoBack bookmark. ---

# Locate that paragraph by its text instead of assuming it is the last one.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*CONFLICTO DE VERSIONES*") {
        $target = $cand
    }
}
if ($target -eq $null) {
    $target = $d.Paragraphs.Item($d.Paragraphs.Count)
}

$r = $target.Range

# Preserve whatever attributes (w14:paraId, rsid*, ...) Word already put
# on that paragraph's opening tag so the rewrite stays as close as
# possible to the original markup.
$openTag = "<w:p>"
$full = $r.WordOpenXML
if ($full -match "<w:p [^>]*>") {
    $openTag = $matches[0]
}

$xml = $openTag + '<w:pPr><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>CONFLICTO DE VERSIONES</w:t></w:r></w:p><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$r.InsertXML($xml)
